# "Ironed out a few bugs, created 'visible' column."
#
# 1) Rename sheet "burp" -> "Process"
# 2) Input sheet: insert a new column G ("Visible") of booleans, fix D4 (9 -> 8)
# 3) Process sheet: minor column width tweak
# 4) Output sheet: VLOOKUP formula now points at Process (no longer wrapped in
#    redundant parens) - values recompute from the corrected Input!D4
# 5) Selections / active sheet follow what the author left the workbook on

$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item("Input")
$wsProcess = $wb.Worksheets.Item("burp")
$wsOutput  = $wb.Worksheets.Item("Output")

# --- rename sheet ---------------------------------------------------------
$wsProcess.Name = "Process"

# --- Input sheet: fix bad width value -------------------------------------
$wsInput.Range("D4").Value = 8

# --- Input sheet: insert new "Visible" column before the old G (Options) --
$wsInput.Columns.Item(7).Insert()

$wsInput.Range("G1").Value = "Visible"
$wsInput.Range("G2").Value = $true
$wsInput.Range("G3").Value = $true
$wsInput.Range("G4").Value = $true
$wsInput.Range("G5").Value = $true

# --- Process sheet: column A width tweak ----------------------------------
$wsProcess.Columns.Item(1).ColumnWidth = 11.7

# --- Output sheet: simplify / repoint the weight VLOOKUP ------------------
$wsOutput.Range("C5").Formula = "=VLOOKUP(Input!D5,Process!A2:B5,2,FALSE)*C3/1000"

# --- selections / active sheet --------------------------------------------
[void]$wsProcess.Activate()
[void]$wsProcess.Range("B6").Select()

[void]$wsOutput.Activate()
[void]$wsOutput.Range("C6").Select()

[void]$wsInput.Activate()
[void]$wsInput.Range("G6").Select()
